$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2 = 2; 3 = 2; 4 = 3; 5 = 12; 6 = 5; 7 = 4; 8 = 7; 9 = 2; 10 = 5;
    11 = 3; 12 = 10; 13 = 4; 14 = 5; 15 = 7; 16 = 4; 17 = 4; 18 = 6;
    19 = 7; 20 = 3; 21 = 4; 22 = 3; 23 = 7; 24 = 5; 25 = 7; 26 = 7;
    27 = 11; 28 = 3; 29 = 3; 30 = 8; 31 = 12; 32 = 5; 33 = 3; 34 = 2
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
